$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# New task row added at the bottom of the list ("revisar reportes de venta - devolucion")
$ws.Range("A43").Value = "revisar reportes de venta - devolucion"

# Fill in the missing progress percentage for row 11 (50%), same percentage
# style already used by the rest of column C.
$ws.Range("C11").Value = 0.5
$ws.Range("C11").NumberFormat = $ws.Range("C2").NumberFormat

# Scroll the view down to show the newly added row and leave the selection
# on the next empty row below it.
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("A44").Select()
